# Apply IFRS data corrections: replace absolute financial figures with
# corrected (smaller-scale) values for FY2014-FY2018, and wipe the
# stale projected rows (FY2019E-FY2021E) which no longer have data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 102
$ws.Range("E2").Value = -16
$ws.Range("F2").Value = -16
$ws.Range("G2").Value = -30
$ws.Range("H2").Value = -35
$ws.Range("I2").Value = -35
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1098
$ws.Range("L2").Value = 48
$ws.Range("M2").Value = 1050
$ws.Range("N2").Value = 1024
$ws.Range("O2").Value = 26
$ws.Range("P2").Value = 287
$ws.Range("Q2").Value = -18
$ws.Range("R2").Value = 108
$ws.Range("S2").Value = 779
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = -20
$ws.Range("W2").Value = -15.91
$ws.Range("X2").Value = -33.8
$ws.Range("Y2").Value = -5.34
$ws.Range("Z2").Value = -4.56
$ws.Range("AA2").Value = 4.56
$ws.Range("AB2").Value = 257.05
$ws.Range("AC2").Value = -130
$ws.Range("AD2").Value = -29.6
$ws.Range("AE2").Value = 1823
$ws.Range("AF2").Value = 2.11
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 56162276
$ws.Range("V2").ClearContents()

# --- Row 3 ---
$ws.Range("D3").Value = 292
$ws.Range("E3").Value = -70
$ws.Range("F3").Value = -70
$ws.Range("G3").Value = -48
$ws.Range("H3").Value = -45
$ws.Range("I3").Value = -36
$ws.Range("J3").Value = -9
$ws.Range("K3").Value = 1165
$ws.Range("L3").Value = 124
$ws.Range("M3").Value = 1040
$ws.Range("N3").Value = 1031
$ws.Range("O3").Value = 10
$ws.Range("P3").Value = 294
$ws.Range("Q3").Value = -59
$ws.Range("R3").Value = -780
$ws.Range("S3").Value = 76
$ws.Range("T3").Value = 37
$ws.Range("U3").Value = -96
$ws.Range("V3").Value = 10
$ws.Range("W3").Value = -23.91
$ws.Range("X3").Value = -15.47
$ws.Range("Y3").Value = -3.49
$ws.Range("Z3").Value = -4
$ws.Range("AA3").Value = 11.94
$ws.Range("AB3").Value = 252.54
$ws.Range("AC3").Value = -63
$ws.Range("AD3").Value = -46.83
$ws.Range("AE3").Value = 1792
$ws.Range("AF3").Value = 1.64
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 57510854

# --- Row 4 ---
$ws.Range("D4").Value = 704
$ws.Range("E4").Value = -68
$ws.Range("F4").Value = -68
$ws.Range("G4").Value = -48
$ws.Range("H4").Value = -54
$ws.Range("I4").Value = -41
$ws.Range("J4").Value = -13
$ws.Range("K4").Value = 1230
$ws.Range("L4").Value = 157
$ws.Range("M4").Value = 1073
$ws.Range("N4").Value = 1032
$ws.Range("O4").Value = 40
$ws.Range("P4").Value = 297
$ws.Range("Q4").Value = -46
$ws.Range("R4").Value = -134
$ws.Range("S4").Value = 93
$ws.Range("T4").Value = 47
$ws.Range("U4").Value = -92
$ws.Range("V4").Value = 34
$ws.Range("W4").Value = -9.609999999999999
$ws.Range("X4").Value = -7.74
$ws.Range("Y4").Value = -4.01
$ws.Range("Z4").Value = -4.55
$ws.Range("AA4").Value = 14.68
$ws.Range("AB4").Value = 239.85
$ws.Range("AC4").Value = -71
$ws.Range("AD4").Value = -29.68
$ws.Range("AE4").Value = 1749
$ws.Range("AF4").Value = 1.2
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 58170704

# --- Row 5 ---
$ws.Range("D5").Value = 725
$ws.Range("E5").Value = -71
$ws.Range("F5").Value = -60
$ws.Range("G5").Value = -89
$ws.Range("H5").Value = -100
$ws.Range("I5").Value = -82
$ws.Range("J5").Value = -18
$ws.Range("K5").Value = 1360
$ws.Range("L5").Value = 219
$ws.Range("M5").Value = 1141
$ws.Range("N5").Value = 958
$ws.Range("O5").Value = 183
$ws.Range("P5").Value = 297
$ws.Range("Q5").Value = -51
$ws.Range("R5").Value = -165
$ws.Range("S5").Value = 216
$ws.Range("T5").Value = 22
$ws.Range("U5").Value = -73
$ws.Range("V5").Value = 85
$ws.Range("W5").Value = -9.779999999999999
$ws.Range("X5").Value = -13.78
$ws.Range("Y5").Value = -8.199999999999999
$ws.Range("Z5").Value = -7.71
$ws.Range("AA5").Value = 19.17
$ws.Range("AB5").Value = 215.15
$ws.Range("AC5").Value = -138
$ws.Range("AD5").Value = -14.47
$ws.Range("AE5").Value = 1623
$ws.Range("AF5").Value = 1.23
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 58170704
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()

# --- Row 6 ---
$ws.Range("D6").Value = 1024
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = -13
$ws.Range("H6").Value = -14
$ws.Range("I6").Value = -1
$ws.Range("K6").Value = 1468
$ws.Range("L6").Value = 312
$ws.Range("M6").Value = 1155
$ws.Range("N6").Value = 961
$ws.Range("P6").Value = 297
$ws.Range("Q6").Value = 91
$ws.Range("R6").Value = -56
$ws.Range("S6").Value = 11
$ws.Range("T6").Value = 21
$ws.Range("U6").Value = 70
$ws.Range("V6").Value = 107
$ws.Range("W6").Value = 0.82
$ws.Range("X6").Value = -1.4
$ws.Range("Y6").Value = -0.15
$ws.Range("Z6").Value = -1.01
$ws.Range("AA6").Value = 27.05
$ws.Range("AB6").Value = 215.8
$ws.Range("AC6").Value = -2
$ws.Range("AD6").Value = -846.5700000000001
$ws.Range("AE6").Value = 1628
$ws.Range("AF6").Value = 1.24
$ws.Range("AJ6").Value = 58170704
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").ClearContents()

# --- Rows 7-9: clear all figures except the row label columns (A-C) ---
$ws.Range("D7:AJ9").ClearContents()
